# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-26 19:19:36
#
# The "Recorded By" column (G) lists the users who recorded/edited a session,
# separated by ", ". This edit re-syncs the ordering of those names to match
# the upstream source-of-truth report (the latest editor is listed first).
#
# Concretely, for every row in column G whose value is one of the known
# mis-ordered combinations, swap the names into the corrected order:
#   "System, backup@backdoor.com, system"  -> "System, system, backup@backdoor.com"
#   "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com"  -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $current = $cell.Value()

    if ($current -eq "System, backup@backdoor.com, system") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
    elseif ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}

$wb.Save()
